$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.566.69"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.814.30"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "229.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.582"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.48%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "34.97"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +7.15%  "
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("E10").Value = "  +0.15%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0954"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "2.075.11"
$ws.Range("E12").Value = "  +0.74%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "11.24"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").Value = "1.812.66"
$ws.Range("E14").Value = "  +0.89%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.647"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "34.548.87"
$ws.Range("E16").Value = "  -0.11%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.88%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "69.26"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "0.0₃0800"
$ws.Range("E19").Value = "  -0.66%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "246.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.55%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.49"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -0.46%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "174.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("E25").Value = "  +2.33%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.93"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +8.51%  "
$ws.Range("E27").Value = "  +1.50%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.120"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("E29").Value = "  -0.16%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0533"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +0.20%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.686"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "1.397.56"
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -0.54%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "83.84"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("E45").Value = "  +3.98%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0513"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.71%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "6.00"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("D48").Value = "1.974.76"
$ws.Range("E48").Value = "  +0.70%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "105.06"
$c.Style = "Normal"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0130"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
